$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 153, shifting existing rows 153-239 down to 154-240.
$ws.Rows.Item(153).Insert()

# Populate the newly inserted row 153 with its data.
$ws.Cells.Item(153, 1).Value = 11
$ws.Cells.Item(153, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(153, 3).Value = "Bíobío"
$ws.Cells.Item(153, 4).Value = 45176
$ws.Cells.Item(153, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(153, 5).Value = 8
$ws.Cells.Item(153, 6).Value = 100112043
$ws.Cells.Item(153, 7).Value = "Pepino ensalada"
$ws.Cells.Item(153, 8).Value = "Sin especificar"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 100
$ws.Cells.Item(153, 11).Value = 13000
$ws.Cells.Item(153, 12).Value = 14000
$ws.Cells.Item(153, 13).Value = 13500
$ws.Cells.Item(153, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(153, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(153, 16).Value = 225
$ws.Cells.Item(153, 17).Value = 60
$ws.Cells.Item(153, 18).Value = "Hortaliza"
